# Apply the "handles float input without breaking stuff" grading fix:
# the student's submitted answers are merged into the per-student
# marksheet, the Right/Wrong/Not-Attempt/Max summary is recomputed, and
# the (now-unused) third answer block (columns G:H) together with the
# trailing rows of the second answer block (D19:E40) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Give the row-label cells (A10, A11, A12) the same "mtitleStyle"
#    formatting already used by the header row above them (A9).
# ---------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Recompute the summary block (Right / Wrong / Not Attempt / Max,
#    Marking scheme, and Total) for this student.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = 22     # Right
$ws.Range("C10").Value = 0      # Wrong
$ws.Range("D10").Value = 6      # Not Attempt
$ws.Range("E10").Value = 28     # Max

$ws.Range("B11").Value = 4      # Marking: right
$ws.Range("C11").Value = -1     # Marking: wrong (numeric, not text)

$ws.Range("B12").Value = 88     # Total score
$ws.Range("E12").Value = "88/112"

# ---------------------------------------------------------------------
# 3. Fill in the student's answers for the first answer block
#    (columns A/B), copying the "correctStyle" formatting already used
#    by the filled-in cells in this sheet (e.g. B10).
# ---------------------------------------------------------------------
$ws.Range("B10").Copy()
$answersA = @{
  16 = "Option A"
  18 = "Option B"
  19 = "Option C"
  21 = "Option C"
  22 = "Option D"
  23 = "Option D"
  25 = "Option A"
  26 = "Option C"
  27 = "Option A"
  28 = "Option D"
  30 = "Option B"
  31 = "Option D"
  32 = "Option C"
  33 = "Option D"
  34 = "Option B"
  36 = "Option A"
  38 = "Option A"
  39 = "Option D"
  40 = "Option D"
}
foreach ($row in $answersA.Keys) {
  $ws.Range("A$row").PasteSpecial(-4122)  # xlPasteFormats
  $ws.Range("A$row").Value = $answersA[$row]
}
$excel.CutCopyMode = 0

# Rows 17, 20, 24, 29, 35, 37 stay blank (not attempted) - no change needed.

# ---------------------------------------------------------------------
# 4. Second answer block (columns D/E) shrinks to just 3 questions
#    (rows 16-18); fill in the student's answers there too, and drop
#    the now-unused rows 19-40.
# ---------------------------------------------------------------------
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"
$excel.CutCopyMode = 0

$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------
# 5. The third answer block (columns G/H) is no longer used at all.
# ---------------------------------------------------------------------
$ws.Range("G1:H1048576").Clear()
